$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.531.87"
$ws.Range("E2").Value = "  +0.61%  "

Set-TextValue $ws.Range("D3") "3.596.94"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue $ws.Range("D5") "609.73"
$ws.Range("E5").Value = "  +0.69%  "

Set-TextValue $ws.Range("D6") "148.99"
$ws.Range("E6").Value = "  +3.11%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -0.43%  "

Set-TextValue $ws.Range("D9") "8.06"
$ws.Range("E9").Value = "  +2.45%  "

$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("E11").Value = "  +0.89%  "

Set-TextValue $ws.Range("D12") "4.210.46"
$ws.Range("E12").Value = "  +1.06%  "

$ws.Range("E13").Value = "  +0.81%  "

Set-TextValue $ws.Range("D14") "29.89"
$ws.Range("E14").Value = "  -0.31%  "

Set-TextValue $ws.Range("D15") "3.607.68"
$ws.Range("E15").Value = "  +1.41%  "

Set-TextValue $ws.Range("D16") "66.604.12"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("E17").Value = "  +0.80%  "

Set-TextValue $ws.Range("D18") "11.56"
$ws.Range("E18").Value = "  +2.17%  "

$ws.Range("E19").Value = "  +3.46%  "

Set-TextValue $ws.Range("D20") "15.13"
$ws.Range("E20").Value = "  +2.10%  "

Set-TextValue $ws.Range("D21") "427.89"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("E22").Value = "  +1.32%  "

Set-TextValue $ws.Range("D23") "78.82"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("E24").Value = "  -0.07%  "

Set-TextValue $ws.Range("D25") "0.0000122"
$ws.Range("E25").Value = "  +3.16%  "

$ws.Range("E26").Value = "  +4.77%  "

Set-TextValue $ws.Range("D27") "9.49"
$ws.Range("E27").Value = "  +4.41%  "

$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("E29").Value = "  -0.01%  "

Set-TextValue $ws.Range("D30") "1.47"
$ws.Range("E30").Value = "  +0.73%  "

Set-TextValue $ws.Range("D31") "3.596.27"
$ws.Range("E31").Value = "  +1.13%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D32") "0.157"
$ws.Range("E32").Value = "  +3.70%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D33") "25.46"
$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  -0.01%  "

Set-TextValue $ws.Range("D36") "5.66"
$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("E37").Value = "  -2.31%  "

Set-TextValue $ws.Range("D38") "177.20"
$ws.Range("E38").Value = "  +0.82%  "

Set-TextValue $ws.Range("D39") "0.0858"
$ws.Range("E39").Value = "  +1.22%  "

Set-TextValue $ws.Range("D41") "0.899"
$ws.Range("E41").Value = "  +0.61%  "

$ws.Range("E42").Value = "  -2.04%  "

$ws.Range("E43").Value = "  +9.97%  "

Set-TextValue $ws.Range("D44") "0.999"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D45") "25.07"
$ws.Range("E45").Value = "  -2.29%  "

$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D46") "1.18"
$ws.Range("E46").Value = "  -2.08%  "

Set-TextValue $ws.Range("D47") "24.02"
$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("E48").Value = "  +1.50%  "

Set-TextValue $ws.Range("D49") "0.956"
$ws.Range("E49").Value = "  +1.59%  "

Set-TextValue $ws.Range("D50") "2.428.54"
$ws.Range("E50").Value = "  +5.50%  "

$ws.Range("E51").Value = "  -0.42%  "

Write-Output "Update complete"